$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 1).Value = "2025-10-28 01:45:11"
}
